$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from diff. Force text number format to preserve values as strings
# (matching the original inlineStr storage) rather than being auto-converted to numbers/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '309.92'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.60%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.20'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.87%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.124'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.11%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07838'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.97%'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '8.281'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.44%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.875'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-1.44%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.965'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '4.31%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9256'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.11%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1177'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.69%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1894'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.76%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08880'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-4.14%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03313'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.45%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09597'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.10%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001376'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.84%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006192'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '6.32%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.389'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-4.02%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.405'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.10%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3458'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.40%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.384'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '21.30%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1293'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.73%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2405'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-6.97%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.61%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-3.79%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004284'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.19%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001399'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '7.93%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002898'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02155'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.60%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05006'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.40%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007567'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.54%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1356'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.37%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008475'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.94%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002012'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.16%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.007866'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-8.61%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006577'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-1.53%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000749'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.21%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003291'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '14.19%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.001442'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '20.50%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002098'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.21%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0001999'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.21%'
